# Update the date header and the multiplication problems/answers in the table.
# Each Find.Execute call targets a unique original text value and replaces it
# with its updated value. Calls are issued in document order so that a value
# that becomes identical to an earlier (already-updated) original value does
# not get erroneously re-matched.

$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-08-09 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-10 Saturday", 2) | Out-Null
$d.Content.Find.Execute("451×2=902", $true, $false, $false, $false, $false, $true, 1, $false, "652×4=2608", 2) | Out-Null
$d.Content.Find.Execute("464×9=4176", $true, $false, $false, $false, $false, $true, 1, $false, "860×8=6880", 2) | Out-Null
$d.Content.Find.Execute("415×8=3320", $true, $false, $false, $false, $false, $true, 1, $false, "694×4=2776", 2) | Out-Null
$d.Content.Find.Execute("625×2=1250", $true, $false, $false, $false, $false, $true, 1, $false, "220×8=1760", 2) | Out-Null
$d.Content.Find.Execute("479×7=3353", $true, $false, $false, $false, $false, $true, 1, $false, "864×5=4320", 2) | Out-Null
$d.Content.Find.Execute("344×9=3096", $true, $false, $false, $false, $false, $true, 1, $false, "511×8=4088", 2) | Out-Null
$d.Content.Find.Execute("443×4=1772", $true, $false, $false, $false, $false, $true, 1, $false, "290×4=1160", 2) | Out-Null
$d.Content.Find.Execute("533×8=4264", $true, $false, $false, $false, $false, $true, 1, $false, "147×7=1029", 2) | Out-Null
$d.Content.Find.Execute("693×5=3465", $true, $false, $false, $false, $false, $true, 1, $false, "300×8=2400", 2) | Out-Null
$d.Content.Find.Execute("532×6=3192", $true, $false, $false, $false, $false, $true, 1, $false, "523×4=2092", 2) | Out-Null
$d.Content.Find.Execute("107×5=535", $true, $false, $false, $false, $false, $true, 1, $false, "446×7=3122", 2) | Out-Null
$d.Content.Find.Execute("470×9=4230", $true, $false, $false, $false, $false, $true, 1, $false, "963×8=7704", 2) | Out-Null
$d.Content.Find.Execute("656×7=4592", $true, $false, $false, $false, $false, $true, 1, $false, "775×7=5425", 2) | Out-Null
$d.Content.Find.Execute("264×4=1056", $true, $false, $false, $false, $false, $true, 1, $false, "299×7=2093", 2) | Out-Null
$d.Content.Find.Execute("797×3=2391", $true, $false, $false, $false, $false, $true, 1, $false, "107×8=856", 2) | Out-Null
$d.Content.Find.Execute("660×7=4620", $true, $false, $false, $false, $false, $true, 1, $false, "285×7=1995", 2) | Out-Null
$d.Content.Find.Execute("831×8=6648", $true, $false, $false, $false, $false, $true, 1, $false, "434×9=3906", 2) | Out-Null
$d.Content.Find.Execute("193×6=1158", $true, $false, $false, $false, $false, $true, 1, $false, "262×2=524", 2) | Out-Null
$d.Content.Find.Execute("140×5=700", $true, $false, $false, $false, $false, $true, 1, $false, "451×2=902", 2) | Out-Null
$d.Content.Find.Execute("615×8=4920", $true, $false, $false, $false, $false, $true, 1, $false, "283×4=1132", 2) | Out-Null
$d.Content.Find.Execute("802×6=4812", $true, $false, $false, $false, $false, $true, 1, $false, "425×6=2550", 2) | Out-Null
$d.Content.Find.Execute("122×5=610", $true, $false, $false, $false, $false, $true, 1, $false, "578×3=1734", 2) | Out-Null
$d.Content.Find.Execute("899×3=2697", $true, $false, $false, $false, $false, $true, 1, $false, "404×6=2424", 2) | Out-Null
$d.Content.Find.Execute("240×7=1680", $true, $false, $false, $false, $false, $true, 1, $false, "576×4=2304", 2) | Out-Null
$d.Content.Find.Execute("640×2=1280", $true, $false, $false, $false, $false, $true, 1, $false, "729×5=3645", 2) | Out-Null
